$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.219.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.88%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.931.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.75%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.34%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'332.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.44%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4736"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -4.94%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.4059"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.66%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'52.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.13%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.08466"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -8.25%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.052"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.30%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.52%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.984.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.52%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.536"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.34%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.124"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.09%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'90.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001063"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06591"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.89%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'18.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.54%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.785"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.82%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'28.230.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.93%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.86%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.290"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.14%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.132.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.37%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'154.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.86%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'20.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.37%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.166"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.87%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.782"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -9.53%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'123.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.30%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.9848"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.98%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09627"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.27%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.454"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.77%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.590"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.14%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.640"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.53%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'9.179"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.63%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02322"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -4.68%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.06180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.95%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.244"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -6.11%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.6195"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.09%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'11.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.73%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.06%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.1905"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.08%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.316"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.52%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5901"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'12.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.61%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.047"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -7.19%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.477"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.33%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06806"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'EOS"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'1.087"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.58%  "
$ws.Range("E51").Style = "Normal"
